# Increment the "Data" column (F) dates by one day for rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @{
    2  = "28/08/2023"
    3  = "28/08/2023"
    4  = "28/08/2023"
    5  = "29/08/2023"
    6  = "29/08/2023"
    7  = "29/08/2023"
    8  = "29/08/2023"
    9  = "29/08/2023"
    10 = "28/08/2023"
    11 = "28/08/2023"
    12 = "28/08/2023"
    13 = "29/08/2023"
    14 = "29/08/2023"
    15 = "29/08/2023"
    16 = "29/08/2023"
    17 = "29/08/2023"
    18 = "28/08/2023"
    19 = "28/08/2023"
    20 = "28/08/2023"
    21 = "29/08/2023"
    22 = "29/08/2023"
    23 = "29/08/2023"
    24 = "29/08/2023"
    25 = "29/08/2023"
}

foreach ($row in $newDates.Keys) {
    $ws.Cells.Item($row, 6).Value = $newDates[$row]
}
